$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as exact text, preserving the cell's original style
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "245.00"
Set-TextValue "E2" "-0.67%"
Set-TextValue "D3" "27.38"
Set-TextValue "E3" "4.78%"
Set-TextValue "D4" "5.118"
Set-TextValue "E4" "0.67%"
Set-TextValue "E5" "1.68%"
Set-TextValue "D6" "6.498"
Set-TextValue "E6" "0.49%"
Set-TextValue "D7" "0.8195"
Set-TextValue "E7" "0.74%"
Set-TextValue "D8" "0.8528"
Set-TextValue "E8" "0.96%"
Set-TextValue "B9" "WazirX"
Set-TextValue "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1331"
Set-TextValue "E9" "-0.31%"
Set-TextValue "B10" "MandalaExchangeToken"
Set-TextValue "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.06954"
Set-TextValue "E10" "0.03%"
Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.02875"
Set-TextValue "E11" "2.04%"
Set-TextValue "B12" "BitMartToken"
Set-TextValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.09394"
Set-TextValue "E12" "0.07%"
Set-TextValue "B13" "BitForexToken"
Set-TextValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001515"
Set-TextValue "E13" "0.05%"
Set-TextValue "B14" "CoinExToken"
Set-TextValue "C14" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D14" "0.04026"
Set-TextValue "E14" "-13.46%"
Set-TextValue "B15" "One"
Set-TextValue "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006030"
Set-TextValue "E15" "-93.94%"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006215"
Set-TextValue "E16" "0.00%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.510"
Set-TextValue "E17" "-2.69%"
Set-TextValue "B18" "GateToken"
Set-TextValue "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D18" "3.010"
Set-TextValue "E18" "-0.32%"
Set-TextValue "B19" "BTSEToken"
Set-TextValue "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D19" "2.319"
Set-TextValue "E19" "12.82%"
Set-TextValue "B20" "BitpandaEcosystemToken"
Set-TextValue "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D20" "0.3158"
Set-TextValue "E20" "1.48%"
Set-TextValue "E21" "1.67%"
Set-TextValue "E22" "-1.54%"
Set-TextValue "D23" "3.555"
Set-TextValue "E23" "-5.51%"
Set-TextValue "E24" "-0.11%"
Set-TextValue "D25" "0.001214"
Set-TextValue "E25" "-2.48%"
Set-TextValue "D26" "0.004477"
Set-TextValue "E26" "-1.64%"
Set-TextValue "E27" "22.81%"
Set-TextValue "E28" "-27.51%"
Set-TextValue "E40" "1.78%"
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.005971"
Set-TextValue "E41" "-3.49%"
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1058"
Set-TextValue "E42" "0.53%"
Set-TextValue "D43" "0.002300"
Set-TextValue "E43" "-11.19%"
Set-TextValue "D44" "0.009711"
Set-TextValue "E44" "20.37%"
Set-TextValue "D45" "0.00005117"
Set-TextValue "E45" "-5.00%"
Set-TextValue "E46" "-0.09%"
Set-TextValue "E47" "-30.42%"
Set-TextValue "D48" "0.002524"
Set-TextValue "E48" "4.21%"
Set-TextValue "E49" "-0.09%"
Set-TextValue "E50" "-0.09%"
